$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C3"   = -11.425
    "D3"   = -7.264999999999999
    "A12"  = -21.669
    "C14"  = -12.586
    "D20"  = -7.755
    "D25"  = -8.061
    "C26"  = -12.868
    "A27"  = -21.908
    "D30"  = -7.206
    "C31"  = -13.05
    "A32"  = -21.702
    "C35"  = -12.589
    "A36"  = -20.178
    "C37"  = -13.329
    "A38"  = -19.741
    "D44"  = -7.746
    "C45"  = -12.441
    "A46"  = -21.828
    "D47"  = -7.375
    "C52"  = -11.363
    "A54"  = -22.15
    "A55"  = -22.278
    "A56"  = -22.061
    "C57"  = -13.829
    "D58"  = -8.062000000000001
    "A67"  = -21.517
    "A69"  = -21.637
    "A72"  = -21.481
    "D78"  = -7.803
    "C81"  = -13.192
    "A83"  = -21.636
    "C83"  = -12.809
    "D84"  = -8.294
    "A86"  = -22.257
    "D89"  = -6.962000000000001
    "A91"  = -21.584
    "D91"  = -6.931999999999999
    "D92"  = -6.837000000000001
    "A93"  = -21.547
    "D96"  = -7.328
    "A99"  = -20.037
    "C100" = -12.638
    "C102" = -13.45
    "D102" = -7.76
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
